$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell E8: "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active cell/selection as the last edited cell
$ws.Range("E8").Select()
